$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = '$ 537,11'
$ws.Range("F3").Value = '$ 1.652,81'
$ws.Range("F4").Value = '$ 909,05'
$ws.Range("F7").Value = '$ 1.570,16'
$ws.Range("F8").Value = '$ 561,90'
$ws.Range("F10").Value = '$ 429,67'
$ws.Range("F11").Value = '$ 1.652,81'
$ws.Range("F12").Value = '$ 272,64'
$ws.Range("F13").Value = '$ 272,64'
$ws.Range("F14").Value = '$ 437,93'
$ws.Range("F15").Value = '$ 495,79'
$ws.Range("F16").Value = '$ 909,00'
$ws.Range("F17").Value = '$ 2.190,00'
$ws.Range("F18").Value = '$ 2.190,00'
$ws.Range("F19").Value = '$ 2.190,00'
$ws.Range("F20").Value = '$ 2.892,48'
$ws.Range("F21").Value = '$ 991,65'
$ws.Range("F22").Value = '$ 1.570,16'
$ws.Range("F23").Value = '$ 991,65'
$ws.Range("F24").Value = '$ 950,33'
$ws.Range("F28").Value = '$ 2.975,12'
$ws.Range("F29").Value = '$ 1.487,52'
$ws.Range("F30").Value = '$ 2.644,54'
$ws.Range("F31").Value = '$ 363,55'
$ws.Range("F32").Value = '$ 437,93'
$ws.Range("F33").Value = '$ 437,93'
$ws.Range("F34").Value = '$ 363,55'
$ws.Range("F35").Value = '$ 743,72'
$ws.Range("F36").Value = '$ 1.652,81'
$ws.Range("F37").Value = '$ 1.652,81'
$ws.Range("F38").Value = '$ 1.429,77'
$ws.Range("F39").Value = '$ 1.239,73'
$ws.Range("F40").Value = '$ 1.330,23'
$ws.Range("F41").Value = '$ 437,93'
$ws.Range("F42").Value = '$ 1.156,94'
$ws.Range("F43").Value = '$ 454,46'
$ws.Range("F44").Value = '$ 1.156,94'
$ws.Range("F45").Value = '$ 1.016,44'
$ws.Range("F46").Value = '$ 437,93'
$ws.Range("F47").Value = '$ 495,79'
$ws.Range("F48").Value = '$ 909,00'
$ws.Range("F51").Value = '$ 1.074,30'
$ws.Range("F52").Value = '$ 661,07'
$ws.Range("F53").Value = '$ 2.661,07'
$ws.Range("F54").Value = '$ 1.652,81'
$ws.Range("F55").Value = '$ 661,07'
$ws.Range("F56").Value = '$ 702,40'
$ws.Range("F57").Value = '$ 454,46'
$ws.Range("F58").Value = '$ 842,89'
$ws.Range("F59").Value = '$ 1.983,38'
$ws.Range("F60").Value = '$ 578,43'
$ws.Range("F61").Value = '$ 578,43'
$ws.Range("F62").Value = '$ 991,65'
$ws.Range("F65").Value = '$ 429,66'
$ws.Range("F69").Value = '$ 479,26'
$ws.Range("F70").Value = '$ 479,26'
$ws.Range("F72").Value = '$ 1.239,90'
$ws.Range("F73").Value = '$ 999,90'
$ws.Range("F74").Value = '$ 1.299,90'
$ws.Range("F75").Value = '$ 1.049,90'
$ws.Range("F76").Value = '$ 1.074,29'
$ws.Range("F77").Value = '$ 702,39'
$ws.Range("F78").Value = '$ 1.818,10'
$ws.Range("F79").Value = '$ 2.272,64'
$ws.Range("F80").Value = '$ 1.322,23'
$ws.Range("F81").Value = '$ 1.032,97'
$ws.Range("F82").Value = '$ 1.032,97'
$ws.Range("F83").Value = '$ 1.032,97'
$ws.Range("F84").Value = '$ 1.032,97'
$ws.Range("F85").Value = '$ 1.032,97'
$ws.Range("F86").Value = '$ 2.479,26'
$ws.Range("F87").Value = '$ 2.768,51'
$ws.Range("F89").Value = '$ 2.768,51'
$ws.Range("F90").Value = '$ 2.768,51'
$ws.Range("F91").Value = '$ 991,65'
$ws.Range("F92").Value = '$ 2.479,26'
